$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values are written as literal text (matching the
# source inlineStr cells), so numeric-looking strings are protected from
# Excel's automatic number coercion via a temporary Text format.
$updates = @(
    @("D2", "56.582.51"),
    @("E2", "  +2.61%  "),
    @("D3", "2.321.16"),
    @("E3", "  +1.14%  "),
    @("D4", "1.01"),
    @("E4", "  +0.55%  "),
    @("D5", "520.46"),
    @("E5", "  +2.56%  "),
    @("D6", "134.50"),
    @("E6", "  +3.48%  "),
    @("D7", "0.997"),
    @("E7", "  +0.14%  "),
    @("E8", "  +1.34%  "),
    @("D9", "2.342.19"),
    @("E9", "  +1.10%  "),
    @("D10", "0.104"),
    @("E10", "  +5.60%  "),
    @("E11", "  -0.73%  "),
    @("D12", "5.24"),
    @("E12", "  +4.02%  "),
    @("D13", "0.342"),
    @("E13", "  +0.24%  "),
    @("D14", "23.78"),
    @("E14", "  -0.30%  "),
    @("D15", "2.760.39"),
    @("E15", "  +2.15%  "),
    @("D16", "56.706.51"),
    @("E16", "  +3.02%  "),
    @("E17", "  +1.94%  "),
    @("D18", "2.329.29"),
    @("E18", "  +0.55%  "),
    @("D19", "10.45"),
    @("E19", "  -3.02%  "),
    @("D20", "4.21"),
    @("E20", "  +0.40%  "),
    @("D21", "322.94"),
    @("E21", "  +3.44%  "),
    @("E22", "  -1.04%  "),
    @("D23", "0.999"),
    @("E23", "  -0.02%  "),
    @("D24", "60.76"),
    @("E24", "  +0.81%  "),
    @("E25", "  +8.58%  "),
    @("D26", "0.998"),
    @("E26", "  +0.54%  "),
    @("D27", "7.89"),
    @("E27", "  +4.72%  "),
    @("E28", "  +13.55%  "),
    @("D29", "0.0₃0747"),
    @("E29", "  +5.25%  "),
    @("B30", "Monero"),
    @("C30", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D30", "169.82"),
    @("E30", "  -1.69%  "),
    @("B31", "PancakeSwap"),
    @("C31", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"),
    @("D31", "1.72"),
    @("E31", "  +4.95%  "),
    @("D32", "6.16"),
    @("E32", "  +0.05%  "),
    @("D33", "18.33"),
    @("E34", "  +0.08%  "),
    @("E35", "  -0.21%  "),
    @("D36", "1.24"),
    @("E36", "  +0.60%  "),
    @("D37", "0.924"),
    @("E37", "  +0.19%  "),
    @("E38", "  +3.73%  "),
    @("D39", "1.55"),
    @("E39", "  +7.24%  "),
    @("D40", "37.94"),
    @("E40", "  +3.23%  "),
    @("E41", "  +0.42%  "),
    @("D42", "3.58"),
    @("E42", "  +4.29%  "),
    @("D43", "137.87"),
    @("E43", "  +1.62%  "),
    @("D44", "5.19"),
    @("E44", "  +5.52%  "),
    @("D45", "275.56"),
    @("E45", "  +6.17%  "),
    @("D46", "0.0932"),
    @("E46", "  +2.10%  "),
    @("D47", "0.0504"),
    @("E47", "  -0.15%  "),
    @("D48", "0.562"),
    @("E48", "  +1.93%  "),
    @("E49", "  +3.07%  "),
    @("D50", "17.88"),
    @("E50", "  +7.15%  "),
    @("E51", "  +0.35%  ")
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    # Force text storage so digit/period-only strings (e.g. "1.01") are not
    # reinterpreted as numbers; ClearFormats keeps the cell style untouched
    # (matching the source, which carries no explicit style on these cells).
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.ClearFormats()
}
